# Updating filtered feeds from workflow
# Inserts a new row (genomeweb.com version of the "circulating tumor cell
# assay" article) as row 85, pushing the former row 85 (360dx.com version of
# the "Guardant Health liquid biopsy" article) down to row 86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift row 85 downward to make room for the new row ---------------
$ws.Rows.Item(85).Insert()

# --- 2. Populate the newly-inserted row 85 (genomeweb.com link, paired
#        with the keywords/title the workflow assigned it) ----------------
$ws.Cells.Item(85, 1).Value = "https://www.genomeweb.com/cancer/circulating-tumor-cell-assay-finds-best-responders-amgens-small-cell-lung-cancer-drug"
$ws.Cells.Item(85, 2).Value = "CDx"
$ws.Cells.Item(85, 3).Value = "Guardant Health Liquid Biopsy Nabs FDA Approval as CDx for Pfizer Colorectal Cancer Drug Combo"

# --- 2b. Re-assert row 86 (the row pushed down from the old row 85) with
#         the keywords/title values the workflow produced for it ---------
$ws.Cells.Item(86, 2).Value = "companion diagnostic"
$ws.Cells.Item(86, 3).Value = "Circulating Tumor Cell Assay Finds Best Responders to Amgen's Small Cell Lung Cancer Drug"

# --- 3. Rebuild the hyperlinks collection in row order --------------------
# (Row-insert does not renumber the worksheet's stored hyperlink refs, so the
# safest route is to clear them all and re-add in the correct, final order.)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2, 1), "https://www.360dx.com/regulatory-news-fda-approvals/beckman-coulter-siemens-healthineers-abbott-others-gain-510k")
$ws.Hyperlinks.Add($ws.Cells.Item(3, 1), "https://www.360dx.com/cancer/agilent-gains-ivdr-certification-expanded-use-cdx-assay-keytruda")
$ws.Hyperlinks.Add($ws.Cells.Item(4, 1), "https://www.genomeweb.com/cancer/entrogen-colorectal-cancer-ras-mutation-detection-test-nabs-cms-coverage")
$ws.Hyperlinks.Add($ws.Cells.Item(5, 1), "https://www.360dx.com/cancer/entrogen-colorectal-cancer-ras-mutation-detection-test-nabs-cms-coverage")
$ws.Hyperlinks.Add($ws.Cells.Item(6, 1), "https://www.fiercebiotech.com/medtech/roche-receives-fda-breakthrough-label-ai-powered-lung-cancer-companion-diagnostic-test")
$ws.Hyperlinks.Add($ws.Cells.Item(7, 1), "https://www.360dx.com/cancer/fda-approves-roche-lung-cancer-cdx-assay-use-abbvies-emrelis")
$ws.Hyperlinks.Add($ws.Cells.Item(8, 1), "https://www.genomeweb.com/cancer/tempus-ai-verastem-partner-develop-cdx-assay-ovarian-cancer-combination-treatment")
$ws.Hyperlinks.Add($ws.Cells.Item(9, 1), "https://www.360dx.com/cancer/tempus-ai-verastem-partner-develop-cdx-assay-ovarian-cancer-combination-treatment")
$ws.Hyperlinks.Add($ws.Cells.Item(10, 1), "https://www.biocentury.com/article/656014/odac-votes-for-darzalex-in-smoldering-multiple-myeloma")
$ws.Hyperlinks.Add($ws.Cells.Item(11, 1), "https://www.genomeweb.com/cancer/illumina-offer-pillar-biosciences-cdx-assay-expanded-partnership")
$ws.Hyperlinks.Add($ws.Cells.Item(12, 1), "https://www.360dx.com/cancer/illumina-offer-pillar-biosciences-cdx-assay-expanded-partnership")
$ws.Hyperlinks.Add($ws.Cells.Item(13, 1), "https://www.genomeweb.com/cancer/biocartis-gains-ivdr-class-c-certification-lung-cancer-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(14, 1), "https://www.360dx.com/cancer/biocartis-gains-ivdr-class-c-certification-lung-cancer-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(15, 1), "https://www.fiercebiotech.com/medtech/asco-guardant-blood-test-guides-breast-cancer-therapy-switches-extending-survival-az-backed")
$ws.Hyperlinks.Add($ws.Cells.Item(16, 1), "https://www.sciencedaily.com/releases/2024/12/241210163409.htm")
$ws.Hyperlinks.Add($ws.Cells.Item(17, 1), "https://www.sciencedaily.com/releases/2024/01/240114202019.htm")
$ws.Hyperlinks.Add($ws.Cells.Item(18, 1), "https://www.sciencedaily.com/releases/2019/10/191028104212.htm")
$ws.Hyperlinks.Add($ws.Cells.Item(19, 1), "https://www.sciencedaily.com/releases/2018/04/180426141507.htm")
$ws.Hyperlinks.Add($ws.Cells.Item(20, 1), "https://www.fiercebiotech.com/biotech/cullinan-pens-700m-pact-bcma-bispecific-pair-another-autoimmune-t-cell-engager")
$ws.Hyperlinks.Add($ws.Cells.Item(21, 1), "https://www.360dx.com/immunoassays/randox-laboratories-gains-fda-de-novo-clearance-hemophilia-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(22, 1), "https://www.fiercebiotech.com/medtech/cancer-tester-caris-life-sciences-go-public-400m-nasdaq-ipo")
$ws.Hyperlinks.Add($ws.Cells.Item(23, 1), "https://www.genomeweb.com/cancer/qiagen-incyte-partner-develop-companion-diagnostics-calr-mutant-other-blood-cancers")
$ws.Hyperlinks.Add($ws.Cells.Item(24, 1), "https://www.360dx.com/cancer/qiagen-incyte-partner-develop-companion-diagnostics-calr-mutant-other-blood-cancers")
$ws.Hyperlinks.Add($ws.Cells.Item(25, 1), "https://www.biocentury.com/article/656269/illumina-s-buy-of-somalogic-could-be-tipping-point-for-multiomics-deals-report")
$ws.Hyperlinks.Add($ws.Cells.Item(26, 1), "https://www.genomeweb.com/cancer/metastx-secures-500k-nci-grant-develop-metastatic-prostate-cancer-cdx-test")
$ws.Hyperlinks.Add($ws.Cells.Item(27, 1), "https://www.360dx.com/cancer/metastx-secures-500k-nci-grant-develop-metastatic-prostate-cancer-cdx-test")
$ws.Hyperlinks.Add($ws.Cells.Item(28, 1), "https://www.genomeweb.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu")
$ws.Hyperlinks.Add($ws.Cells.Item(29, 1), "https://www.360dx.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu")
$ws.Hyperlinks.Add($ws.Cells.Item(30, 1), "https://www.genomeweb.com/sequencing/thermo-fisher-ngs-assay-gets-fda-ok-cdx-dizals-zegfrovy-and-solid-tumor-profiling")
$ws.Hyperlinks.Add($ws.Cells.Item(31, 1), "https://www.360dx.com/sequencing/thermo-fisher-ngs-assay-gets-fda-ok-cdx-dizals-zegfrovy-and-solid-tumor-profiling")
$ws.Hyperlinks.Add($ws.Cells.Item(32, 1), "https://www.genomeweb.com/sequencing/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(33, 1), "https://www.360dx.com/sequencing/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(34, 1), "https://www.genomeweb.com/cancer/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(35, 1), "https://www.360dx.com/cancer/thermo-fisher-nabs-fda-approval-ngs-based-cdx-boehringer-ingelheim-lung-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(36, 1), "https://www.360dx.com/cancer/fda-approves-agilent-cdx-assay-use-bristol-myers-squibb-immunotherapies")
$ws.Hyperlinks.Add($ws.Cells.Item(37, 1), "https://www.genomeweb.com/cancer/caris-life-sciences-eyes-wider-clinical-adoption-tests-new-validation-data")
$ws.Hyperlinks.Add($ws.Cells.Item(38, 1), "https://www.360dx.com/cancer/caris-life-sciences-eyes-wider-clinical-adoption-tests-new-validation-data")
$ws.Hyperlinks.Add($ws.Cells.Item(39, 1), "https://www.360dx.com/cancer/agilent-secures-ivdr-class-c-certification-colorectal-cancer-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(40, 1), "https://www.360dx.com/cancer/roche-nabs-ce-ivdr-marking-her2-cdx-assay-breast-biliary-tract-cancer")
$ws.Hyperlinks.Add($ws.Cells.Item(41, 1), "https://www.genomeweb.com/cancer/biocartis-gets-fda-approval-msi-companion-diagnostic")
$ws.Hyperlinks.Add($ws.Cells.Item(42, 1), "https://www.360dx.com/cancer/biocartis-gets-fda-approval-msi-companion-diagnostic")
$ws.Hyperlinks.Add($ws.Cells.Item(43, 1), "https://www.genomeweb.com/cancer/amoydx-lung-cancer-panel-nabs-japanese-approval-cdx-non-small-cell-lung-cancer-drug-ibtrozi")
$ws.Hyperlinks.Add($ws.Cells.Item(44, 1), "https://www.360dx.com/cancer/amoydx-lung-cancer-panel-nabs-japanese-approval-cdx-non-small-cell-lung-cancer-drug-ibtrozi")
$ws.Hyperlinks.Add($ws.Cells.Item(45, 1), "https://www.medpagetoday.com/meetingcoverage/ims/117536")
$ws.Hyperlinks.Add($ws.Cells.Item(46, 1), "https://www.medpagetoday.com/meetingcoverage/ims/117570")
$ws.Hyperlinks.Add($ws.Cells.Item(47, 1), "https://www.genomeweb.com/cancer/ogt-assay-nabs-fda-authorization-cdx-syndaxs-revuforj-acute-leukemia")
$ws.Hyperlinks.Add($ws.Cells.Item(48, 1), "https://www.360dx.com/cancer/agilent-technologies-lunit-partner-cancer-cdx-development")
$ws.Hyperlinks.Add($ws.Cells.Item(49, 1), "https://www.360dx.com/cancer/ogt-assay-nabs-fda-authorization-cdx-syndaxs-revuforj-acute-leukemia")
$ws.Hyperlinks.Add($ws.Cells.Item(50, 1), "https://www.genomeweb.com/cancer/myriad-genetics-sophia-genetics-collaborate-cancer-liquid-biopsy-companion-diagnostic")
$ws.Hyperlinks.Add($ws.Cells.Item(51, 1), "https://www.360dx.com/cancer/myriad-genetics-sophia-genetics-collaborate-cancer-liquid-biopsy-companion-diagnostic")
$ws.Hyperlinks.Add($ws.Cells.Item(52, 1), "https://www.genomeweb.com/cancer/burning-rocks-sequencing-test-approved-japan-cdx-truqap-breast-cancer")
$ws.Hyperlinks.Add($ws.Cells.Item(53, 1), "https://www.360dx.com/cancer/burning-rocks-sequencing-test-approved-japan-cdx-truqap-breast-cancer")
$ws.Hyperlinks.Add($ws.Cells.Item(54, 1), "https://www.genomeweb.com/companion-diagnostics/biocartis-positioning-idylla-system-sample-answer-oncology-cdx-tests")
$ws.Hyperlinks.Add($ws.Cells.Item(55, 1), "https://www.360dx.com/companion-diagnostics/biocartis-positioning-idylla-system-sample-answer-oncology-cdx-tests")
$ws.Hyperlinks.Add($ws.Cells.Item(56, 1), "https://www.genomeweb.com/cancer/guardant-health-blood-test-gets-fda-ok-cdx-eli-lilly-breast-cancer-drug-inluriyo")
$ws.Hyperlinks.Add($ws.Cells.Item(57, 1), "https://www.360dx.com/cancer/guardant-health-blood-test-gets-fda-ok-cdx-eli-lilly-breast-cancer-drug-inluriyo")
$ws.Hyperlinks.Add($ws.Cells.Item(58, 1), "https://www.genomeweb.com/companion-diagnostics/celebrating-10th-anniversary-its-pd-l1-cdx-agilent-sets-sights-next-gen")
$ws.Hyperlinks.Add($ws.Cells.Item(59, 1), "https://www.genomeweb.com/cancer/geneseeq-nabs-china-nmpa-approval-pan-solid-tumor-test-cdx-roches-rozlytrek")
$ws.Hyperlinks.Add($ws.Cells.Item(60, 1), "https://www.360dx.com/cancer/geneseeq-nabs-china-nmpa-approval-pan-solid-tumor-test-cdx-roches-rozlytrek")
$ws.Hyperlinks.Add($ws.Cells.Item(61, 1), "https://www.genomeweb.com/cancer/promega-receives-fda-approval-oncomate-msi-dx-analysis-system-endometrial-cancer-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(62, 1), "https://www.360dx.com/cancer/promega-receives-fda-approval-oncomate-msi-dx-analysis-system-endometrial-cancer-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(63, 1), "https://www.fiercebiotech.com/medtech/promega-msi-tech-wins-fda-approval-companion-diagnostic-keytruda-lenvima-combo-uterine")
$ws.Hyperlinks.Add($ws.Cells.Item(64, 1), "https://www.genomeweb.com/regulatory-news-fda-approvals/thermo-fisher-receives-fda-approval-ngs-based-cdx-bayer-lung-cancer")
$ws.Hyperlinks.Add($ws.Cells.Item(65, 1), "https://www.360dx.com/regulatory-news-fda-approvals/thermo-fisher-receives-fda-approval-ngs-based-cdx-bayer-lung-cancer")
$ws.Hyperlinks.Add($ws.Cells.Item(66, 1), "https://www.genomeweb.com/cancer/fda-proposes-reclassification-companion-diagnostic-tests")
$ws.Hyperlinks.Add($ws.Cells.Item(67, 1), "https://www.360dx.com/cancer/fda-proposes-reclassification-companion-diagnostic-tests")
$ws.Hyperlinks.Add($ws.Cells.Item(68, 1), "https://www.360dx.com/business-news/top-five-articles-360dx-last-week-fda-aims-reclassify-cdx-tests-tempus-buying-oneome")
$ws.Hyperlinks.Add($ws.Cells.Item(69, 1), "https://www.genomeweb.com/cancer/fda-proposal-reclassify-cdx-assays-may-broaden-opportunities-dx-manufacturers-experts-say")
$ws.Hyperlinks.Add($ws.Cells.Item(70, 1), "https://www.360dx.com/cancer/fda-proposal-reclassify-cdx-assays-may-broaden-opportunities-dx-manufacturers-experts-say")
$ws.Hyperlinks.Add($ws.Cells.Item(71, 1), "https://www.360dx.com/business-news/top-five-articles-360dx-last-week-roche-vaginitis-test-fda-reclassification-cdx")
$ws.Hyperlinks.Add($ws.Cells.Item(72, 1), "https://www.genomeweb.com/regulatory-news-fda-approvals/fda-expands-labels-roches-her2-cdx-assays-id-breast-cancer-patients")
$ws.Hyperlinks.Add($ws.Cells.Item(73, 1), "https://www.360dx.com/regulatory-news-fda-approvals/fda-expands-labels-roches-her2-cdx-assays-id-breast-cancer-patients")
$ws.Hyperlinks.Add($ws.Cells.Item(74, 1), "https://www.genomeweb.com/cancer/guardant-health-blood-test-receives-japanese-approval-cdx-eli-lilly-breast-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(75, 1), "https://www.360dx.com/cancer/guardant-health-blood-test-receives-japanese-approval-cdx-eli-lilly-breast-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(76, 1), "https://www.genomeweb.com/cancer/chinese-mdx-precision-medicine-firm-geneplus-files-ipo-hong-kong-stock-exchange")
$ws.Hyperlinks.Add($ws.Cells.Item(77, 1), "https://www.360dx.com/cancer/chinese-mdx-precision-medicine-firm-geneplus-files-ipo-hong-kong-stock-exchange")
$ws.Hyperlinks.Add($ws.Cells.Item(78, 1), "https://www.genomeweb.com/cancer/precede-raises-84m-integrate-liquid-biopsy-test-alongside-next-gen-precision-cancer-drugs")
$ws.Hyperlinks.Add($ws.Cells.Item(79, 1), "https://www.360dx.com/cancer/precede-raises-84m-integrate-liquid-biopsy-test-alongside-next-gen-precision-cancer-drugs")
$ws.Hyperlinks.Add($ws.Cells.Item(80, 1), "https://www.genomeweb.com/cancer/guardant-merck-partner-develop-and-commercialize-companion-diagnostics")
$ws.Hyperlinks.Add($ws.Cells.Item(81, 1), "https://www.360dx.com/cancer/guardant-merck-partner-develop-and-commercialize-companion-diagnostics")
$ws.Hyperlinks.Add($ws.Cells.Item(82, 1), "https://www.fiercebiotech.com/medtech/guardant-health-pens-merck-co-cancer-collab-next-gen-tests-and-companion-diagnostics")
$ws.Hyperlinks.Add($ws.Cells.Item(83, 1), "https://www.genomeweb.com/cancer/guardant-health-liquid-biopsy-nabs-fda-approval-cdx-pfizer-colorectal-cancer-drug-combo")
$ws.Hyperlinks.Add($ws.Cells.Item(84, 1), "https://www.360dx.com/cancer/circulating-tumor-cell-assay-finds-best-responders-amgens-small-cell-lung-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(85, 1), "https://www.genomeweb.com/cancer/circulating-tumor-cell-assay-finds-best-responders-amgens-small-cell-lung-cancer-drug")
$ws.Hyperlinks.Add($ws.Cells.Item(86, 1), "https://www.360dx.com/cancer/guardant-health-liquid-biopsy-nabs-fda-approval-cdx-pfizer-colorectal-cancer-drug-combo")
